$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.62 = 36258.62 pesos`n✅ 36258.62 pesos = 8.61 = 950.13 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 116
$wsTasas.Range("O10").Value = 4206
$wsTasas.Range("N12").Value = 4210
$wsTasas.Range("O12").Value = 110.32
